$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 465.45456
$ws.Range("I6").Value = 184
$ws.Range("J6").Value = 700
$ws.Range("K6").Value = 552
$ws.Range("L6").Value = 2100
$ws.Range("M6").Value = -440
$ws.Range("N6").Value = -2324
$ws.Range("H137").Value = 4374098.5
$ws.Range("I137").Value = 6075960.5
$ws.Range("J137").Value = 2502050
$ws.Range("K137").Value = 18227881.5
$ws.Range("L137").Value = 7506150
$ws.Range("M137").Value = -18225331.5
$ws.Range("N137").Value = -7511250
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1603.909
$ws.Range("I2").Value = 2009
$ws.Range("J2").Value = 895
$ws.Range("K2").Value = 2009
$ws.Range("L2").Value = 895
$ws.Range("M2").Value = -1896
$ws.Range("N2").Value = -1121
$ws.Range("H32").Value = 1606447.5
$ws.Range("I32").Value = 1942148.1
$ws.Range("J32").Value = 2544.2222
$ws.Range("K32").Value = 1942148.1
$ws.Range("L32").Value = 2544.2222
$ws.Range("M32").Value = -1941861.1
$ws.Range("N32").Value = -3118.2222
$ws.Range("H45").Value = 1741.8125
$ws.Range("I45").Value = 1634.2142
$ws.Range("J45").Value = 2495
$ws.Range("K45").Value = 1634.2142
$ws.Range("L45").Value = 2495
$ws.Range("M45").Value = -1257.2142
$ws.Range("N45").Value = -3249
$ws.Range("H61").Value = 479976.75
$ws.Range("I61").Value = 386512.28
$ws.Range("J61").Value = 631856.5600000001
$ws.Range("K61").Value = 386512.28
$ws.Range("L61").Value = 631856.5600000001
$ws.Range("M61").Value = -386300.28
$ws.Range("N61").Value = -632280.5600000001
$ws.Range("H110").Value = 2163.6428
$ws.Range("I110").Value = 2279
$ws.Range("J110").Value = 1956
$ws.Range("K110").Value = 2279
$ws.Range("L110").Value = 1956
$ws.Range("M110").Value = -234
$ws.Range("N110").Value = -6046
$ws.Range("H116").Value = 1603.909
$ws.Range("I116").Value = 2009
$ws.Range("J116").Value = 895
$ws.Range("K116").Value = 2009
$ws.Range("L116").Value = 895
$ws.Range("M116").Value = 285
$ws.Range("N116").Value = -5483
$ws.Range("H122").Value = 2255.0715
$ws.Range("I122").Value = 1359.45
$ws.Range("K122").Value = 4078.35
$ws.Range("M122").Value = -1628.35
$ws.Range("H132").Value = 23394.49
$ws.Range("I132").Value = 42954.418
$ws.Range("K132").Value = 128863.254
$ws.Range("M132").Value = -126333.254
$ws.Range("H136").Value = 479976.75
$ws.Range("I136").Value = 386512.28
$ws.Range("J136").Value = 631856.5600000001
$ws.Range("K136").Value = 1159536.84
$ws.Range("L136").Value = 1895569.68
$ws.Range("M136").Value = -1156986.84
$ws.Range("N136").Value = -1900669.68
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1603.909
$ws.Range("I3").Value = 2009
$ws.Range("J3").Value = 895
$ws.Range("K3").Value = 2009
$ws.Range("L3").Value = 895
$ws.Range("M3").Value = -1895
$ws.Range("N3").Value = -1123
$ws.Range("H20").Value = 1005.75
$ws.Range("I20").Value = 907.05
$ws.Range("J20").Value = 1129.125
$ws.Range("K20").Value = 907.05
$ws.Range("L20").Value = 1129.125
$ws.Range("M20").Value = -660.05
$ws.Range("N20").Value = -1623.125
$ws.Range("H105").Value = 1228.375
$ws.Range("I105").Value = 864.2857
$ws.Range("J105").Value = 3777
$ws.Range("K105").Value = 864.2857
$ws.Range("L105").Value = 3777
$ws.Range("M105").Value = 882.7143
$ws.Range("N105").Value = -7271
$ws.Range("H132").Value = 28000
$ws.Range("J132").Value = 28000
$ws.Range("L132").Value = 28000
$ws.Range("N132").Value = -38120
$ws.Range("H134").Value = 3686.1333
$ws.Range("I134").Value = 3730.1538
$ws.Range("J134").Value = 3400
$ws.Range("K134").Value = 11190.4614
$ws.Range("L134").Value = 10200
$ws.Range("M134").Value = -8655.4614
$ws.Range("N134").Value = -15270
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1050
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1050
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1050
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1624
$ws.Range("H94").Value = 7029.643
$ws.Range("I94").Value = 1335
$ws.Range("J94").Value = 17280
$ws.Range("K94").Value = 1335
$ws.Range("L94").Value = 17280
$ws.Range("M94").Value = -884
$ws.Range("N94").Value = -18182
$ws.Range("H99").Value = 49186.57
$ws.Range("I99").Value = 56863.777
$ws.Range("J99").Value = 3123.3333
$ws.Range("K99").Value = 56863.777
$ws.Range("L99").Value = 3123.3333
$ws.Range("M99").Value = -55365.777
$ws.Range("N99").Value = -6119.3333
$ws.Range("H111").Value = 39000
$ws.Range("J111").Value = 39000
$ws.Range("L111").Value = 39000
$ws.Range("N111").Value = -47180
$ws.Range("H113").Value = 1050
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1050
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1050
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5390
$ws.Range("H126").Value = 49186.57
$ws.Range("I126").Value = 56863.777
$ws.Range("J126").Value = 3123.3333
$ws.Range("K126").Value = 170591.331
$ws.Range("L126").Value = 9369.999899999999
$ws.Range("M126").Value = -168121.331
$ws.Range("N126").Value = -14309.9999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 473.08334
$ws.Range("I7").Value = 219.75
$ws.Range("J7").Value = 599.75
$ws.Range("K7").Value = 659.25
$ws.Range("L7").Value = 1799.25
$ws.Range("M7").Value = -547.25
$ws.Range("N7").Value = -2023.25
$ws.Range("H12").Value = 464.83334
$ws.Range("I12").Value = 233
$ws.Range("J12").Value = 696.6667
$ws.Range("K12").Value = 699
$ws.Range("L12").Value = 2090.0001
$ws.Range("M12").Value = -526
$ws.Range("N12").Value = -2436.0001
$ws.Range("H33").Value = 253.70589
$ws.Range("I33").Value = 165.07143
$ws.Range("J33").Value = 667.3333
$ws.Range("K33").Value = 990.42858
$ws.Range("L33").Value = 4003.9998
$ws.Range("M33").Value = -707.42858
$ws.Range("N33").Value = -4569.9998
$ws.Range("H51").Value = 1258.5454
$ws.Range("J51").Value = 1171.4286
$ws.Range("L51").Value = 3514.2858
$ws.Range("N51").Value = -4434.2858
$ws.Range("H92").Value = 290
$ws.Range("I92").Value = 300
$ws.Range("K92").Value = 900
$ws.Range("M92").Value = 348
$ws.Range("H119").Value = 679.6667
$ws.Range("I119").Value = 679.6667
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 2039.0001
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 2798.9999
$ws.Range("N119").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3804.7058
$ws.Range("I80").Value = 4199.231
$ws.Range("J80").Value = 2522.5
$ws.Range("K80").Value = 4199.231
$ws.Range("L80").Value = 2522.5
$ws.Range("M80").Value = -3201.231
$ws.Range("N80").Value = -4518.5
$ws.Range("H83").Value = 3804.7058
$ws.Range("I83").Value = 4199.231
$ws.Range("J83").Value = 2522.5
$ws.Range("K83").Value = 20996.155
$ws.Range("L83").Value = 12612.5
$ws.Range("M83").Value = -16004.155
$ws.Range("N83").Value = -22596.5
$ws.Range("H102").Value = 9579.583000000001
$ws.Range("I102").Value = 6883.8887
$ws.Range("J102").Value = 17666.666
$ws.Range("K102").Value = 6883.8887
$ws.Range("L102").Value = 17666.666
$ws.Range("M102").Value = -5261.8887
$ws.Range("N102").Value = -20910.666
$ws.Range("H126").Value = 2542
$ws.Range("I126").Value = 2123.75
$ws.Range("K126").Value = 6371.25
$ws.Range("M126").Value = -3901.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2726
$ws.Range("I7").Value = 2634.6667
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2634.6667
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2522.6667
$ws.Range("N7").Value = -3224
$ws.Range("H126").Value = 2726
$ws.Range("I126").Value = 2634.6667
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7904.000100000001
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5434.000100000001
$ws.Range("N126").Value = -13940
